$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")

# Change L10, M10, L11, M11 from the text "NA" to the boolean value FALSE
$ws.Range("L10").Value = $false
$ws.Range("M10").Value = $false
$ws.Range("L11").Value = $false
$ws.Range("M11").Value = $false

# Update the selected cell to Q10
$ws.Range("Q10").Select()
